$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 / Row 7 - swap Starting/Ending SoC (%) values
$ws.Range("B6").Value = 89
$ws.Range("B7").Value = 8

# Row 8
$ws.Range("A8").Value = "Total distance covered (km)"

# Row 9
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"

# Row 10
$ws.Range("A10").Value = "Total SOC consumed(%)"

# Row 12
$ws.Range("A12").Value = "Peak Power(kW)"

# Row 13
$ws.Range("A13").Value = "Average Power(kW)"

# Row 14
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"

# Row 15 - label + sign flip on value
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 0.0376013560653105

# Row 16 / Row 17 - swap Lowest/Highest Cell Voltage (label & value)
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.371
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 2.98

# Row 18
$ws.Range("A18").Value = "Difference in Cell Voltage(V)"

# Row 19
$ws.Range("A19").Value = "Minimum Temperature(C)"

# Row 20
$ws.Range("A20").Value = "Maximum Temperature(C)"

# Row 21 - label + fill in previously empty value
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 13

# Row 22
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"

# Row 23
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"

# Row 24
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"

# Row 25
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"

# Row 26
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"

# Row 27
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

# Row 28 / Row 29 - swap lowest/highest cell temp labels (values stay put)
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("A29").Value = "lowest cell temp(C)"

# Row 30
$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# Row 31 - label & value now represents Battery Voltage
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 53

# Row 32 - label & value now represents Total energy charged
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.689470223333333

# Row 33 - label & value now represents Electricity consumption units
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.0000001213909167768389

# Row 34 - label & value now represents Idling time percentage
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 22.70450751252087

# Row 35 - label & value now represents Time spent in 0-10 km/h
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 2.899276572064552

# Row 36 - label & value now represents Time spent in 10-20 km/h
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 4.774624373956595

# Row 37 - label & value now represents Time spent in 20-30 km/h
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 9.315525876460768

# Row 38 - label & value now represents Time spent in 30-40 km/h
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 30.28380634390651

# Row 39 - label & value now represents Time spent in 40-50 km/h
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 5.779076238174736

# Row 40 - label & value now represents Time spent in 50-60 km/h
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 8.597662771285476

# Row 41 - label & value now represents Time spent in 60-70 km/h
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 11.47746243739566

# Row 42 - label & value now represents Time spent in 70-80 km/h
$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 4.067890929326656

# Row 43 - new row: Time spent in 80-90 km/h
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
